$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# X2 / X3: set numeric value 1
$ws.Range("X2").Value = 1
$ws.Range("X3").Value = 1

# Y2 / Y3: set text value "x"
$ws.Range("Y2").Value = "x"
$ws.Range("Y3").Value = "x"
